# Auto update Excel log
# Appends newly-logged sensor events to the ALERTS, PIR, and mmWave sheets.
#
# Each data row is written as plain text (matching the existing inline-string
# cells already in the sheet). In particular, the Date column ("2026-01-30")
# would otherwise be auto-recognized by Excel as a date literal and stamped
# with a date number format, so we force the cell to Text format before
# assigning the value and then reset its style to "Normal" afterwards so the
# cell keeps the default (unstyled) look used by every other row.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($ws, $startRow, $rows)

    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $startRow + $i
        $data = $rows[$i]

        $dateCell = $ws.Cells.Item($r, 1)
        $dateCell.NumberFormat = "@"
        $dateCell.Value = $data[0]
        $dateCell.Style = "Normal"

        $ws.Cells.Item($r, 2).Value = $data[1]
        $ws.Cells.Item($r, 3).Value = $data[2]
        $ws.Cells.Item($r, 4).Value = $data[3]
        $ws.Cells.Item($r, 5).Value = $data[4]
        $ws.Cells.Item($r, 6).Value = $data[5]
    }
}

# ---------------------------------------------------------------------------
# ALERTS sheet: append rows 6-10
# ---------------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")

$alertsRows = @(
    ,@("2026-01-30","14:14:43","14:00","Living Room","FALL_DETECTED","CRITICAL EMERGENCY")
    ,@("2026-01-30","14:14:47","14:00","Living Room","FALL_DETECTED","CRITICAL EMERGENCY")
    ,@("2026-01-30","14:15:06","14:00","Living Room","FALL_DETECTED","CRITICAL EMERGENCY")
    ,@("2026-01-30","14:15:10","14:00","Living Room","FALL_DETECTED","CRITICAL EMERGENCY")
    ,@("2026-01-30","14:15:44","14:00","Living Room","FALL_DETECTED","CRITICAL EMERGENCY")
)

Add-LogRows $wsAlerts 6 $alertsRows

# ---------------------------------------------------------------------------
# PIR sheet: append row 127
# ---------------------------------------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")

$pirRows = @(
    ,@("2026-01-30","14:05:53","14:00","Living Room","RECOVERY_DETECTION","Inactive")
)

Add-LogRows $wsPir 127 $pirRows

# ---------------------------------------------------------------------------
# mmWave sheet: append rows 50-59
# ---------------------------------------------------------------------------
$wsMmWave = $wb.Worksheets.Item("mmWave")

$mmWaveRows = @(
    ,@("2026-01-30","14:05:37","14:00","Living Room","FALL_DETECTED","EMERGENCY")
    ,@("2026-01-30","14:05:53","14:00","Living Room","PRESENCE_DETECTED","Active")
    ,@("2026-01-30","14:06:03","14:00","Living Room","PRESENCE_DETECTED","Active")
    ,@("2026-01-30","14:14:43","14:00","Living Room","FALL_DETECTED","CRITICAL EMERGENCY")
    ,@("2026-01-30","14:14:47","14:00","Living Room","FALL_DETECTED","CRITICAL EMERGENCY")
    ,@("2026-01-30","14:15:06","14:00","Living Room","FALL_DETECTED","CRITICAL EMERGENCY")
    ,@("2026-01-30","14:15:10","14:00","Living Room","FALL_DETECTED","CRITICAL EMERGENCY")
    ,@("2026-01-30","14:15:11","14:00","Living Room","PRESENCE_DETECTED","Active")
    ,@("2026-01-30","14:15:16","14:00","Living Room","PRESENCE_DETECTED","Active")
    ,@("2026-01-30","14:15:44","14:00","Living Room","FALL_DETECTED","CRITICAL EMERGENCY")
)

Add-LogRows $wsMmWave 50 $mmWaveRows
